$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 69-78 (levels 66-75) - table now ends at level 65 (row 68)
$ws.Rows("69:78").Delete()

# Fix the animus level values for rows 59-68 (levels 56-65) that were
# off-by-a-tiny-float-precision-amount before
$ws.Range("B59:D59").Value = 1285081075117470
$ws.Range("B60:D60").Value = 2088256747065880
$ws.Range("B61:D61").Value = 3257680525422780
$ws.Range("B62:D62").Value = 4940815463557880
$ws.Range("B63:D63").Value = 7340640117286010
$ws.Range("B64:D64").Value = 13121394209648700
$ws.Range("B65:D65").Value = 21709943146873400
$ws.Range("B66:D66").Value = 34270695967564400
$ws.Range("B67:D67").Value = 52414005597451500
$ws.Range("B68:D68").Value = 78358938368189900
